# plne funkcni direct editing
# Applies the "direct editing" changes recorded for TRIMAZKON/saved_addresses_2.xlsx:
#  - updates/retypes rows on "ip_address_list" and adds rows 12-15
#  - updates/retypes rows on "ip_adress_fav_list" and adds rows 3-7
#  - flips a flag on "Settings"
#  - updates the active-cell selection on the two address sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")
$wsSettings = $wb.Worksheets.Item("Settings")

# ---------------------------------------------------------------------------
# Sheet "ip_address_list": columns A,B,C,D,E ; E is a flag (bool or number)
# ---------------------------------------------------------------------------
$sheet1Rows = @(
    @{ R=1;  A="ggfs";       B="192.168.000.000";  C="255.255.255.0"; D="ggiif";                                                            E=1; EType="b" },
    @{ R=2;  A="518_Valeoo"; B="192.168.208.242";  C="255.255.255.0"; D="k";                                                                 E=1; EType="b" },
    @{ R=3;  A="oo";         B="192.168.000.000";  C="255.255.255.0"; D="ooo";                                                               E=1; EType="b" },
    @{ R=4;  A="527_Teihg";  B="10.101.28.17";     C="255.255.255.0"; D="XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28.h";                    E=1; EType="b" },
    @{ R=5;  A="kkkk";       B="10.96.205.17";     C="255.255.255.0"; D="PC:`t10.96.205.175NAS:`t10.96.205.166k`nFH:`t10.96.205.154`n`t10.96.20aa"; E=1; EType="b" },
    @{ R=6;  A="hhggg";      B="192.168.000.000h"; C="255.255.255.0"; D="hhh";                                                               E=0; EType="b" },
    @{ R=7;  A="527_Tei";    B="10.101.28.17";     C="255.255.255.0"; D="XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28.h";                    E=1; EType="b" },
    @{ R=8;  A="se";         B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=0; EType="b" },
    @{ R=9;  A="h";          B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=0; EType="b" },
    @{ R=10; A="gg";         B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=0; EType="b" },
    @{ R=11; A="jjs";        B="192.168.000.000";  C="255.255.255.0"; D="ss";                                                                E=1; EType="b" },
    @{ R=12; A="hhh";        B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=0; EType="n" },
    @{ R=13; A="jjjj";       B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=0; EType="n" },
    @{ R=14; A="aauj";       B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=1; EType="n" },
    @{ R=15; A="ss";         B="192.168.000.000";  C="255.255.255.0"; D=$null;                                                               E=1; EType="n" }
)

foreach ($row in $sheet1Rows) {
    $r = $row.R
    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B
    $ws1.Cells.Item($r, 3).Value = $row.C

    if ($row.D -eq $null) {
        $ws1.Cells.Item($r, 4).ClearContents()
    } else {
        $ws1.Cells.Item($r, 4).Value = $row.D
    }

    if ($row.EType -eq "b") {
        if ($row.E -eq 1) {
            $ws1.Cells.Item($r, 5).Value = $true
        } else {
            $ws1.Cells.Item($r, 5).Value = $false
        }
    } else {
        $ws1.Cells.Item($r, 5).Value = $row.E
    }
}

# ---------------------------------------------------------------------------
# Sheet "ip_adress_fav_list": columns A,B,C,D,E ; E is a flag (bool or number)
# ---------------------------------------------------------------------------
$sheet2Rows = @(
    @{ R=1; A="kkkk";      B="10.96.205.17";    C="255.255.255.0"; D="PC:`t10.96.205.175NAS:`t10.96.205.166k`nFH:`t10.96.205.154`n`t10.96.20aa"; E=1; EType="n" },
    @{ R=2; A="sssssssss"; B="192.168.000.000"; C="255.255.255.0"; D=$null; E=1; EType="n" },
    @{ R=3; A="sssse";     B="192.168.000.000"; C="255.255.255.0"; D=$null; E=1; EType="n" },
    @{ R=4; A="aaaee";     B="192.168.000.000"; C="255.255.255.0"; D=$null; E=1; EType="n" },
    @{ R=5; A="aauj";      B="192.168.000.000"; C="255.255.255.0"; D=$null; E=1; EType="n" },
    @{ R=6; A="jjs";       B="192.168.000.000"; C="255.255.255.0"; D="ss";  E=1; EType="b" },
    @{ R=7; A="ss";        B="192.168.000.000"; C="255.255.255.0"; D=$null; E=1; EType="n" }
)

foreach ($row in $sheet2Rows) {
    $r = $row.R
    $ws2.Cells.Item($r, 1).Value = $row.A
    $ws2.Cells.Item($r, 2).Value = $row.B
    $ws2.Cells.Item($r, 3).Value = $row.C

    if ($row.D -eq $null) {
        $ws2.Cells.Item($r, 4).ClearContents()
    } else {
        $ws2.Cells.Item($r, 4).Value = $row.D
    }

    if ($row.EType -eq "b") {
        if ($row.E -eq 1) {
            $ws2.Cells.Item($r, 5).Value = $true
        } else {
            $ws2.Cells.Item($r, 5).Value = $false
        }
    } else {
        $ws2.Cells.Item($r, 5).Value = $row.E
    }
}

$ws2.Range("A8:C26").Select()

# ---------------------------------------------------------------------------
# Sheet "Settings": flip B6 from 0 to 1
# ---------------------------------------------------------------------------
$wsSettings.Cells.Item(6, 2).Value = 1

# ---------------------------------------------------------------------------
# Restore "ip_address_list" as the active sheet/selection (selected last so it
# stays the active tab, matching tabSelected="1" on that sheet).
# ---------------------------------------------------------------------------
$ws1.Range("C23").Select()
